# Fix a bug in DownView.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B9 previously held the numeric hours value 1; it now holds the text "1+2"
# (two separate work sessions summed together).
$ws.Range("B9").Value = "1+2"

# A brand-new day of work (row 10) is logged: 2 hours spent loading images
# into the Photo Wall feature.
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "Load images into Photo Wall"

# The active selection moves from C10 to B9 (the cell that was just edited).
$ws.Range("B9").Select()
